# Apply leve-profit recalculations per the commit diff.
# Each block targets one (sheet, row); sets changed H..N cells, clears any that were removed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (ALC)
$ws.Range("H32").Value = 3282.3333
$ws.Range("I32").Value = 4600.6665
$ws.Range("K32").Value = 4600.6665
$ws.Range("M32").Value = -4274.6665

# Row 33 (ALC)
$ws.Range("H33").Value = 720411.5600000001
$ws.Range("I33").Value = 1015535.5
$ws.Range("J33").Value = 3682
$ws.Range("K33").Value = 1015535.5
$ws.Range("L33").Value = 3682
$ws.Range("M33").Value = -1015306.5
$ws.Range("N33").Value = -4140

# Row 112 (ALC)
$ws.Range("H112").Value = 2187.6667
$ws.Range("J112").Value = 2187.6667
$ws.Range("L112").Value = 6563.000100000001
$ws.Range("N112").Value = -8779.000100000001

# Row 116 (ALC)
$ws.Range("H116").Value = 4989.6875
$ws.Range("I116").Value = 4962.3335
$ws.Range("J116").Value = 4996
$ws.Range("K116").Value = 4962.3335
$ws.Range("L116").Value = 4996
$ws.Range("M116").Value = -1520.3335
$ws.Range("N116").Value = -11880

# Row 137 (ALC)
$ws.Range("H137").Value = 49317.383
$ws.Range("I137").Value = 2166.6155
$ws.Range("J137").Value = 125937.375
$ws.Range("K137").Value = 6499.8465
$ws.Range("L137").Value = 377812.125
$ws.Range("M137").Value = -3949.8465
$ws.Range("N137").Value = -382912.125

$ws = $wb.Worksheets.Item("ARM")
# Row 15 (ARM)
$ws.Range("H15").Value = 4295
$ws.Range("J15").Value = 4295
$ws.Range("L15").Value = 4295
$ws.Range("N15").Value = -4995

# Row 37 (ARM)
$ws.Range("H37").Value = 20678
$ws.Range("J37").Value = 32000
$ws.Range("L37").Value = 32000
$ws.Range("N37").Value = -32546

# Row 125 (ARM)
$ws.Range("H125").Value = 100715
$ws.Range("J125").Value = 100715
$ws.Range("L125").Value = 100715
$ws.Range("N125").Value = -110555

# Row 132 (ARM)
$ws.Range("H132").Value = 15417.125
$ws.Range("I132").Value = 16691
$ws.Range("K132").Value = 50073
$ws.Range("M132").Value = -47543

$ws = $wb.Worksheets.Item("BSM")
# Row 35 (BSM)
$ws.Range("H35").Value = 29999.8
$ws.Range("J35").Value = 29999.8
$ws.Range("L35").Value = 29999.8
$ws.Range("N35").Value = -30619.8

# Row 105 (BSM)
$ws.Range("H105").Value = 2503536.5
$ws.Range("I105").Value = 3707612
$ws.Range("K105").Value = 3707612
$ws.Range("M105").Value = -3705865

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 1306.1428
$ws.Range("J16").Value = 1010.2
$ws.Range("L16").Value = 1010.2
$ws.Range("N16").Value = -1584.2

# Row 31 (CRP)
$ws.Range("H31").Value = 1939.9025
$ws.Range("I31").Value = 1544.8334
$ws.Range("J31").Value = 3017.3635
$ws.Range("K31").Value = 1544.8334
$ws.Range("L31").Value = 3017.3635
$ws.Range("M31").Value = -1249.8334
$ws.Range("N31").Value = -3607.3635

# Row 34 (CRP)
$ws.Range("H34").Value = 1939.9025
$ws.Range("I34").Value = 1544.8334
$ws.Range("J34").Value = 3017.3635
$ws.Range("K34").Value = 1544.8334
$ws.Range("L34").Value = 3017.3635
$ws.Range("M34").Value = -1342.8334
$ws.Range("N34").Value = -3421.3635

# Row 99 (CRP)
$ws.Range("H99").Value = 201663.2
$ws.Range("I99").Value = 251500.5
$ws.Range("J99").Value = 2314
$ws.Range("K99").Value = 251500.5
$ws.Range("L99").Value = 2314
$ws.Range("M99").Value = -250002.5
$ws.Range("N99").Value = -5310

# Row 113 (CRP)
$ws.Range("H113").Value = 1306.1428
$ws.Range("J113").Value = 1010.2
$ws.Range("L113").Value = 1010.2
$ws.Range("N113").Value = -5350.2

# Row 122 (CRP)
$ws.Range("H122").Value = 5829.4287
$ws.Range("I122").Value = 5829.4287
$ws.Range("K122").Value = 17488.2861
$ws.Range("M122").Value = -15038.2861

# Row 126 (CRP)
$ws.Range("H126").Value = 201663.2
$ws.Range("I126").Value = 251500.5
$ws.Range("J126").Value = 2314
$ws.Range("K126").Value = 754501.5
$ws.Range("L126").Value = 6942
$ws.Range("M126").Value = -752031.5
$ws.Range("N126").Value = -11882

# Row 132 (CRP)
$ws.Range("H132").Value = 3781.2
$ws.Range("I132").Value = 3582.2
$ws.Range("J132").Value = 3980.2
$ws.Range("K132").Value = 10746.6
$ws.Range("L132").Value = 11940.6
$ws.Range("M132").Value = -8216.599999999999
$ws.Range("N132").Value = -17000.6

$ws = $wb.Worksheets.Item("CUL")
# Row 11 (CUL)
$ws.Range("H11").Value = 67082
$ws.Range("I11").Value = 400.54544
$ws.Range("J11").Value = 250456
$ws.Range("K11").Value = 1201.63632
$ws.Range("L11").Value = 751368
$ws.Range("M11").Value = -1061.63632
$ws.Range("N11").Value = -751648

# Row 48 (CUL)
$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 2000
$ws.Range("K48").Value = 6000
$ws.Range("M48").Value = -5750

# Row 131 (CUL)
$ws.Range("H131").Value = 7708.1113
$ws.Range("I131").Value = 12305.111
$ws.Range("J131").Value = 3111.111
$ws.Range("K131").Value = 36915.333
$ws.Range("L131").Value = 9333.332999999999
$ws.Range("M131").Value = -31875.333
$ws.Range("N131").Value = -19413.333

$ws = $wb.Worksheets.Item("GSM")
# Row 19 (GSM)
$ws.Range("H19").Value = 5491.2856
$ws.Range("I19").Value = 9000
$ws.Range("J19").Value = 4087.8
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 4087.8
$ws.Range("M19").Value = -8712
$ws.Range("N19").Value = -4663.8

# Row 43 (GSM)
$ws.Range("H43").Value = 26699.5
$ws.Range("I43").Value = 3400
$ws.Range("J43").Value = 49999
$ws.Range("K43").Value = 3400
$ws.Range("L43").Value = 49999
$ws.Range("M43").Value = -3249
$ws.Range("N43").Value = -50301

# Row 46 (GSM)
$ws.Range("H46").Value = 19999.777
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 34999.5
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 34999.5
$ws.Range("M46").Value = -7844
$ws.Range("N46").Value = -35311.5

# Row 57 (GSM)
$ws.Range("H57").Value = 25666.5
$ws.Range("J57").Value = 25999.8
$ws.Range("L57").Value = 25999.8
$ws.Range("N57").Value = -27639.8

# Row 126 (GSM)
$ws.Range("H126").Value = 3009.9
$ws.Range("I126").Value = 3166.5
$ws.Range("K126").Value = 9499.5
$ws.Range("M126").Value = -7029.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 975
$ws.Range("I7").Value = 975
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 975
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -863
$ws.Range("N7").ClearContents()

# Row 126 (LTW)
$ws.Range("H126").Value = 975
$ws.Range("I126").Value = 975
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2925
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -455
$ws.Range("N126").ClearContents()

# Row 132 (LTW)
$ws.Range("H132").Value = 7641.5713
$ws.Range("I132").Value = 8971.357
$ws.Range("J132").Value = 4982
$ws.Range("K132").Value = 26914.071
$ws.Range("L132").Value = 14946
$ws.Range("M132").Value = -24384.071
$ws.Range("N132").Value = -20006

# Row 136 (LTW)
$ws.Range("H136").Value = 2928.75
$ws.Range("I136").Value = 1739.1666
$ws.Range("K136").Value = 5217.4998
$ws.Range("M136").Value = -2667.4998

$ws = $wb.Worksheets.Item("WVR")
# Row 18 (WVR)
$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 5000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -4827
$ws.Range("N18").ClearContents()

# Row 48 (WVR)
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# Row 95 (WVR)
$ws.Range("H95").Value = 37867.2
$ws.Range("J95").Value = 37867.2
$ws.Range("L95").Value = 37867.2
$ws.Range("N95").Value = -43359.2

# Row 136 (WVR)
$ws.Range("H136").Value = 952.375
$ws.Range("I136").Value = 952.375
$ws.Range("K136").Value = 2857.125
$ws.Range("M136").Value = -307.125
